# Auto-generated edit script
# Applies updated market-board price / profit figures to the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5567.222
$ws.Range("I62").Value = 2651.25
$ws.Range("J62").Value = 7900
$ws.Range("K62").Value = 2651.25
$ws.Range("L62").Value = 7900
$ws.Range("M62").Value = -2027.25
$ws.Range("N62").Value = -9148
$ws.Range("H65").Value = 5567.222
$ws.Range("I65").Value = 2651.25
$ws.Range("J65").Value = 7900
$ws.Range("K65").Value = 13256.25
$ws.Range("L65").Value = 39500
$ws.Range("M65").Value = -10136.25
$ws.Range("N65").Value = -45740
$ws.Range("H86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 12441.353
$ws.Range("J44").Value = 12441.353
$ws.Range("L44").Value = 12441.353
$ws.Range("N44").Value = -13417.353
$ws.Range("H45").Value = 2657.875
$ws.Range("I45").Value = 2311.5
$ws.Range("K45").Value = 2311.5
$ws.Range("M45").Value = -1934.5
$ws.Range("H61").Value = 2283.6428
$ws.Range("I61").Value = 2283.6428
$ws.Range("K61").Value = 2283.6428
$ws.Range("M61").Value = -2071.6428
$ws.Range("H74").Value = 2749.8572
$ws.Range("I74").Value = 2338.6924
$ws.Range("K74").Value = 2338.6924
$ws.Range("M74").Value = -1464.6924
$ws.Range("H77").Value = 2749.8572
$ws.Range("I77").Value = 2338.6924
$ws.Range("K77").Value = 11693.462
$ws.Range("M77").Value = -7325.462
$ws.Range("H97").Value = 932.7895
$ws.Range("I97").Value = 802.2857
$ws.Range("J97").Value = 1298.2
$ws.Range("K97").Value = 802.2857
$ws.Range("L97").Value = 1298.2
$ws.Range("M97").Value = -306.2857
$ws.Range("N97").Value = -2290.2
$ws.Range("H122").Value = 1985.3077
$ws.Range("I122").Value = 1649.5
$ws.Range("J122").Value = 3104.6667
$ws.Range("K122").Value = 4948.5
$ws.Range("L122").Value = 9314.000100000001
$ws.Range("M122").Value = -2498.5
$ws.Range("N122").Value = -14214.0001
$ws.Range("H132").Value = 2869.077
$ws.Range("I132").Value = 2750
$ws.Range("K132").Value = 8250
$ws.Range("M132").Value = -5720
$ws.Range("H136").Value = 2283.6428
$ws.Range("I136").Value = 2283.6428
$ws.Range("K136").Value = 6850.928400000001
$ws.Range("M136").Value = -4300.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4332
$ws.Range("I86").Value = 4332
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4332
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3209
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 4332
$ws.Range("I89").Value = 4332
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 21660
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -16044
$ws.Range("N89").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 5415.9287
$ws.Range("I134").Value = 1247.9231
$ws.Range("J134").Value = 59600
$ws.Range("K134").Value = 3743.7693
$ws.Range("L134").Value = 178800
$ws.Range("M134").Value = -1208.7693
$ws.Range("N134").Value = -183870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 355.5
$ws.Range("I16").Value = 355.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 355.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -68.5
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 6672.9688
$ws.Range("I31").Value = 4289.1665
$ws.Range("J31").Value = 7223.077
$ws.Range("K31").Value = 4289.1665
$ws.Range("L31").Value = 7223.077
$ws.Range("M31").Value = -3994.1665
$ws.Range("N31").Value = -7813.077
$ws.Range("H34").Value = 6672.9688
$ws.Range("I34").Value = 4289.1665
$ws.Range("J34").Value = 7223.077
$ws.Range("K34").Value = 4289.1665
$ws.Range("L34").Value = 7223.077
$ws.Range("M34").Value = -4087.1665
$ws.Range("N34").Value = -7627.077
$ws.Range("H68").Value = 74382
$ws.Range("J68").Value = 74382
$ws.Range("L68").Value = 74382
$ws.Range("N68").Value = -75880
$ws.Range("H71").Value = 74382
$ws.Range("J71").Value = 74382
$ws.Range("L71").Value = 223146
$ws.Range("N71").Value = -230634
$ws.Range("H99").Value = 2352.3333
$ws.Range("I99").Value = 1522.4286
$ws.Range("K99").Value = 1522.4286
$ws.Range("M99").Value = -24.42859999999996
$ws.Range("H105").Value = 1300
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 567.2308
$ws.Range("I107").Value = 299.2
$ws.Range("J107").Value = 734.75
$ws.Range("K107").Value = 299.2
$ws.Range("L107").Value = 734.75
$ws.Range("M107").Value = 1620.8
$ws.Range("N107").Value = -4574.75
$ws.Range("H113").Value = 355.5
$ws.Range("I113").Value = 355.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 355.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1814.5
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2352.3333
$ws.Range("I126").Value = 1522.4286
$ws.Range("K126").Value = 4567.2858
$ws.Range("M126").Value = -2097.2858
$ws.Range("H132").Value = 1383.5834
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -155060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1430
$ws.Range("I92").Value = 781.3333
$ws.Range("J92").Value = 2264
$ws.Range("K92").Value = 2343.9999
$ws.Range("L92").Value = 6792
$ws.Range("M92").Value = -1095.9999
$ws.Range("N92").Value = -9288
$ws.Range("H107").Value = 338.33334
$ws.Range("J107").Value = 338.33334
$ws.Range("L107").Value = 1015.00002
$ws.Range("N107").Value = -4855.00002
$ws.Range("H123").Value = 330
$ws.Range("I123").Value = 330
$ws.Range("K123").Value = 990
$ws.Range("M123").Value = 1460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 20232382
$ws.Range("J3").Value = 13752625
$ws.Range("L3").Value = 13752625
$ws.Range("N3").Value = -13752857
$ws.Range("H80").Value = 4733.1665
$ws.Range("I80").Value = 4699.75
$ws.Range("J80").Value = 4800
$ws.Range("K80").Value = 4699.75
$ws.Range("L80").Value = 4800
$ws.Range("M80").Value = -3701.75
$ws.Range("N80").Value = -6796
$ws.Range("H83").Value = 4733.1665
$ws.Range("I83").Value = 4699.75
$ws.Range("J83").Value = 4800
$ws.Range("K83").Value = 23498.75
$ws.Range("L83").Value = 24000
$ws.Range("M83").Value = -18506.75
$ws.Range("N83").Value = -33984
$ws.Range("H132").Value = 59399.777
$ws.Range("I132").Value = 70160.39999999999
$ws.Range("K132").Value = 210481.2
$ws.Range("M132").Value = -207951.2
$ws.Range("H140").Value = 118237.75
$ws.Range("I140").Value = 279697
$ws.Range("J140").Value = 95172.14
$ws.Range("K140").Value = 279697
$ws.Range("L140").Value = 95172.14
$ws.Range("M140").Value = -274517
$ws.Range("N140").Value = -105532.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 696.25
$ws.Range("I46").Value = 523.2
$ws.Range("J46").Value = 819.8570999999999
$ws.Range("K46").Value = 523.2
$ws.Range("L46").Value = 819.8570999999999
$ws.Range("M46").Value = -335.2
$ws.Range("N46").Value = -1195.8571
$ws.Range("H61").Value = 4712.7144
$ws.Range("I61").Value = 2247.25
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 2247.25
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -2045.25
$ws.Range("N61").Value = -8404
$ws.Range("H82").Value = 3604.5454
$ws.Range("I82").Value = 1230
$ws.Range("K82").Value = 1230
$ws.Range("M82").Value = -869
$ws.Range("H85").Value = 3604.5454
$ws.Range("I85").Value = 1230
$ws.Range("K85").Value = 1230
$ws.Range("M85").Value = 18
$ws.Range("H113").Value = 4712.7144
$ws.Range("I113").Value = 2247.25
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 2247.25
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -77.25
$ws.Range("N113").Value = -12340
$ws.Range("H122").Value = 1149.75
$ws.Range("I122").Value = 1066.3334
$ws.Range("K122").Value = 3199.0002
$ws.Range("M122").Value = -749.0001999999999
$ws.Range("H136").Value = 4000.8
$ws.Range("I136").Value = 3999.75
$ws.Range("J136").Value = 4005
$ws.Range("K136").Value = 11999.25
$ws.Range("L136").Value = 12015
$ws.Range("M136").Value = -9449.25
$ws.Range("N136").Value = -17115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 12500
$ws.Range("J97").Value = 12500
$ws.Range("L97").Value = 12500
$ws.Range("N97").Value = -14482
$ws.Range("H132").Value = 614.4
$ws.Range("I132").Value = 614.4
$ws.Range("K132").Value = 1843.2
$ws.Range("M132").Value = 686.8000000000002
$ws.Range("H136").Value = 2009.0392
$ws.Range("I136").Value = 1486.5676
$ws.Range("J136").Value = 3389.8572
$ws.Range("K136").Value = 4459.7028
$ws.Range("L136").Value = 10169.5716
$ws.Range("M136").Value = -1909.7028
$ws.Range("N136").Value = -15269.5716
